$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.059713262147824
$ws.Cells.Item(2, 4).Value = 1.061176994324406
$ws.Cells.Item(2, 5).Value = 1.072006602950772
$ws.Cells.Item(2, 6).Value = 1.076397512472445
$ws.Cells.Item(2, 9).Value = 1.044010326904563
$ws.Cells.Item(2, 10).Value = 1.064697865535582
$ws.Cells.Item(2, 11).Value = 1.063901764446436
$ws.Cells.Item(2, 12).Value = 1.074702264085774
$ws.Cells.Item(2, 13).Value = 1.07908155160653
$ws.Cells.Item(2, 14).Value = 1.02518816096133
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.060997101664716
$ws.Cells.Item(3, 4).Value = 1.062148334255206
$ws.Cells.Item(3, 5).Value = 1.07321857082966
$ws.Cells.Item(3, 6).Value = 1.077553285550385
$ws.Cells.Item(3, 9).Value = 1.044275865559087
$ws.Cells.Item(3, 10).Value = 1.065633470943978
$ws.Cells.Item(3, 11).Value = 1.06468727905079
$ws.Cells.Item(3, 12).Value = 1.075729918898881
$ws.Cells.Item(3, 13).Value = 1.080053993895812
$ws.Cells.Item(3, 14).Value = 1.025509223645544
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.061827690667976
$ws.Cells.Item(4, 4).Value = 1.062776532244934
$ws.Cells.Item(4, 5).Value = 1.074002938183728
$ws.Cells.Item(4, 6).Value = 1.078301153011383
$ws.Cells.Item(4, 9).Value = 1.044446219915029
$ws.Cells.Item(4, 10).Value = 1.066238212468257
$ws.Cells.Item(4, 11).Value = 1.065194615896669
$ws.Cells.Item(4, 12).Value = 1.076394449570731
$ws.Cells.Item(4, 13).Value = 1.080682647852486
$ws.Cells.Item(4, 14).Value = 1.025716521119377
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.0621768389115
$ws.Cells.Item(5, 4).Value = 1.063040550099304
$ws.Cells.Item(5, 5).Value = 1.074332722065717
$ws.Cells.Item(5, 6).Value = 1.078615559052529
$ws.Cells.Item(5, 9).Value = 1.044517486057025
$ws.Cells.Item(5, 10).Value = 1.066492289728136
$ws.Cells.Item(5, 11).Value = 1.065407675408987
$ws.Cells.Item(5, 12).Value = 1.076673716578266
$ws.Cells.Item(5, 13).Value = 1.08094679608874
$ws.Cells.Item(5, 14).Value = 1.025803561222187
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.062235460667494
$ws.Cells.Item(6, 4).Value = 1.063084875395952
$ws.Cells.Item(6, 5).Value = 1.074388096378839
$ws.Cells.Item(6, 6).Value = 1.078668349400075
$ws.Cells.Item(6, 9).Value = 1.044529431393274
$ws.Cells.Item(6, 10).Value = 1.066534941305916
$ws.Cells.Item(6, 11).Value = 1.065443435879793
$ws.Cells.Item(6, 12).Value = 1.076720600816112
$ws.Cells.Item(6, 13).Value = 1.080991139690751
$ws.Cells.Item(6, 14).Value = 1.025818169329694
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.061832356124186
$ws.Cells.Item(7, 4).Value = 1.06278006036414
$ws.Cells.Item(7, 5).Value = 1.074007344629261
$ws.Cells.Item(7, 6).Value = 1.078305354108527
$ws.Cells.Item(7, 9).Value = 1.044447173554787
$ws.Cells.Item(7, 10).Value = 1.066241608072585
$ws.Cells.Item(7, 11).Value = 1.06519746369074
$ws.Cells.Item(7, 12).Value = 1.076398181547935
$ws.Cells.Item(7, 13).Value = 1.080686177954378
$ws.Cells.Item(7, 14).Value = 1.025717684576462
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.060147172356862
$ws.Cells.Item(8, 4).Value = 1.061505330846725
$ws.Cells.Item(8, 5).Value = 1.072416164498649
$ws.Cells.Item(8, 6).Value = 1.076788111183606
$ws.Cells.Item(8, 9).Value = 1.04410037088164
$ws.Cells.Item(8, 10).Value = 1.065014194788114
$ws.Cells.Item(8, 11).Value = 1.064167428777691
$ws.Cells.Item(8, 12).Value = 1.075049654881675
$ws.Cells.Item(8, 13).Value = 1.079410314138297
$ws.Cells.Item(8, 14).Value = 1.025296759305034
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.057176476294237
$ws.Cells.Item(9, 4).Value = 1.059256570319524
$ws.Cells.Item(9, 5).Value = 1.069613326174818
$ws.Cells.Item(9, 6).Value = 1.074114511977683
$ws.Cells.Item(9, 9).Value = 1.043478009490856
$ws.Cells.Item(9, 10).Value = 1.06284623582936
$ws.Cells.Item(9, 11).Value = 1.062345102983589
$ws.Cells.Item(9, 12).Value = 1.072670015643965
$ws.Cells.Item(9, 13).Value = 1.077157570685407
$ws.Cells.Item(9, 14).Value = 1.024551561073593
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.055195062554057
$ws.Cells.Item(10, 4).Value = 1.057755635848056
$ws.Cells.Item(10, 5).Value = 1.067745353367342
$ws.Cells.Item(10, 6).Value = 1.072332008028302
$ws.Cells.Item(10, 9).Value = 1.043055513930957
$ws.Cells.Item(10, 10).Value = 1.061397411436417
$ws.Cells.Item(10, 11).Value = 1.061125263407886
$ws.Cells.Item(10, 12).Value = 1.071081236232868
$ws.Cells.Item(10, 13).Value = 1.075652628576138
$ws.Cells.Item(10, 14).Value = 1.024052403347545
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.054336828865312
$ws.Cells.Item(11, 4).Value = 1.057105282017761
$ws.Cells.Item(11, 5).Value = 1.066936615876632
$ws.Cells.Item(11, 6).Value = 1.071560118751765
$ws.Cells.Item(11, 9).Value = 1.042870763279022
$ws.Cells.Item(11, 10).Value = 1.060769199341388
$ws.Cells.Item(11, 11).Value = 1.060595868591142
$ws.Cells.Item(11, 12).Value = 1.070392697783132
$ws.Cells.Item(11, 13).Value = 1.075000214629203
$ws.Cells.Item(11, 14).Value = 1.023835697570722
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.054017999271312
$ws.Cells.Item(12, 4).Value = 1.05686364425463
$ws.Cells.Item(12, 5).Value = 1.066636228712215
$ws.Cells.Item(12, 6).Value = 1.071273394820238
$ws.Cells.Item(12, 9).Value = 1.042801866587724
$ws.Cells.Item(12, 10).Value = 1.060535722087765
$ws.Cells.Item(12, 11).Value = 1.060399046645561
$ws.Cells.Item(12, 12).Value = 1.070136853905021
$ws.Cells.Item(12, 13).Value = 1.074757762509563
$ws.Cells.Item(12, 14).Value = 1.02375511766887
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.05408639131411
$ws.Cells.Item(13, 4).Value = 1.05691547944724
$ws.Cells.Item(13, 5).Value = 1.066700662181847
$ws.Cells.Item(13, 6).Value = 1.071334898554444
$ws.Cells.Item(13, 9).Value = 1.042816657482598
$ws.Cells.Item(13, 10).Value = 1.060585809723374
$ws.Cells.Item(13, 11).Value = 1.060441273840804
$ws.Cells.Item(13, 12).Value = 1.070191737390932
$ws.Cells.Item(13, 13).Value = 1.074809774613474
$ws.Cells.Item(13, 14).Value = 1.023772406218989
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.054310475214434
$ws.Cells.Item(14, 4).Value = 1.057085309547776
$ws.Cells.Item(14, 5).Value = 1.066911785507747
$ws.Cells.Item(14, 6).Value = 1.071536418274694
$ws.Cells.Item(14, 9).Value = 1.042865073810416
$ws.Cells.Item(14, 10).Value = 1.060749902728292
$ws.Cells.Item(14, 11).Value = 1.060579602931137
$ws.Cells.Item(14, 12).Value = 1.070371551496337
$ws.Cells.Item(14, 13).Value = 1.074980175848622
$ws.Cells.Item(14, 14).Value = 1.023829038562696
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.054448534860132
$ws.Cells.Item(15, 4).Value = 1.057189938505707
$ws.Cells.Item(15, 5).Value = 1.067041867341101
$ws.Cells.Item(15, 6).Value = 1.071660579859078
$ws.Cells.Item(15, 9).Value = 1.042894868639149
$ws.Cells.Item(15, 10).Value = 1.060850988388472
$ws.Cells.Item(15, 11).Value = 1.060664808003863
$ws.Cells.Item(15, 12).Value = 1.070482328912228
$ws.Cells.Item(15, 13).Value = 1.075085150180016
$ws.Cells.Item(15, 14).Value = 1.023863920243607
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.055252014780586
$ws.Cells.Item(16, 4).Value = 1.057798788315548
$ws.Cells.Item(16, 5).Value = 1.067799028684591
$ws.Cells.Item(16, 6).Value = 1.072383234494288
$ws.Cells.Item(16, 9).Value = 1.043067737099015
$ws.Cells.Item(16, 10).Value = 1.06143908549248
$ws.Cells.Item(16, 11).Value = 1.061160372274141
$ws.Cells.Item(16, 12).Value = 1.071126919741191
$ws.Cells.Item(16, 13).Value = 1.07569591089473
$ws.Cells.Item(16, 14).Value = 1.024066773402454
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.055755942476019
$ws.Cells.Item(17, 4).Value = 1.05818058512042
$ws.Cells.Item(17, 5).Value = 1.068274003148411
$ws.Cells.Item(17, 6).Value = 1.072836521059848
$ws.Cells.Item(17, 9).Value = 1.043175688667759
$ws.Cells.Item(17, 10).Value = 1.061807751283123
$ws.Cells.Item(17, 11).Value = 1.061470905382678
$ws.Cells.Item(17, 12).Value = 1.071531096160737
$ws.Cells.Item(17, 13).Value = 1.076078819164064
$ws.Cells.Item(17, 14).Value = 1.024193865615323
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.056049849279983
$ws.Cells.Item(18, 4).Value = 1.058403238446843
$ws.Cells.Item(18, 5).Value = 1.068551058193012
$ws.Cells.Item(18, 6).Value = 1.073100910429932
$ws.Cells.Item(18, 9).Value = 1.043238480652122
$ws.Cells.Item(18, 10).Value = 1.062022704799825
$ws.Cells.Item(18, 11).Value = 1.061651918769273
$ws.Cells.Item(18, 12).Value = 1.071766788972989
$ws.Cells.Item(18, 13).Value = 1.076302089475089
$ws.Cells.Item(18, 14).Value = 1.024267941662783
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.056150059638445
$ws.Cells.Item(19, 4).Value = 1.0584791503457
$ws.Cells.Item(19, 5).Value = 1.068645528644372
$ws.Cells.Item(19, 6).Value = 1.073191059578333
$ws.Cells.Item(19, 9).Value = 1.043259861541514
$ws.Cells.Item(19, 10).Value = 1.062095984371641
$ws.Cells.Item(19, 11).Value = 1.061713620138183
$ws.Cells.Item(19, 12).Value = 1.071847144613152
$ws.Cells.Item(19, 13).Value = 1.076378206448772
$ws.Cells.Item(19, 14).Value = 1.024293190419624
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.055701878484563
$ws.Cells.Item(20, 4).Value = 1.05813962632184
$ws.Cells.Item(20, 5).Value = 1.068223041816579
$ws.Cells.Item(20, 6).Value = 1.072787888234568
$ws.Cells.Item(20, 9).Value = 1.043164124523163
$ws.Cells.Item(20, 10).Value = 1.061768205532534
$ws.Cells.Item(20, 11).Value = 1.061437600058359
$ws.Cells.Item(20, 12).Value = 1.071487737710397
$ws.Cells.Item(20, 13).Value = 1.076037744357331
$ws.Cells.Item(20, 14).Value = 1.024180235484575
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.054244489313258
$ws.Cells.Item(21, 4).Value = 1.057035300687692
$ws.Cells.Item(21, 5).Value = 1.066849614582518
$ws.Cells.Item(21, 6).Value = 1.071477076027839
$ws.Cells.Item(21, 9).Value = 1.042850823923585
$ws.Cells.Item(21, 10).Value = 1.060701585068102
$ws.Cells.Item(21, 11).Value = 1.060538873473553
$ws.Cells.Item(21, 12).Value = 1.070318603225128
$ws.Cells.Item(21, 13).Value = 1.07493000016281
$ws.Cells.Item(21, 14).Value = 1.02381236411964
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.053327918343619
$ws.Cells.Item(22, 4).Value = 1.056340576717391
$ws.Cells.Item(22, 5).Value = 1.065986166009173
$ws.Cells.Item(22, 6).Value = 1.0706528582362
$ws.Cells.Item(22, 9).Value = 1.042652264997456
$ws.Cells.Item(22, 10).Value = 1.060030197797019
$ws.Cells.Item(22, 11).Value = 1.059972759676266
$ws.Cells.Item(22, 12).Value = 1.069583000064123
$ws.Cells.Item(22, 13).Value = 1.074232842344747
$ws.Cells.Item(22, 14).Value = 1.023580572542182
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.053813834918355
$ws.Cells.Item(23, 4).Value = 1.056708900531782
$ws.Cells.Item(23, 5).Value = 1.066443889601387
$ws.Cells.Item(23, 6).Value = 1.071089797821514
$ws.Cells.Item(23, 9).Value = 1.042757674260195
$ws.Cells.Item(23, 10).Value = 1.060386185619289
$ws.Cells.Item(23, 11).Value = 1.060272967193823
$ws.Cells.Item(23, 12).Value = 1.069973007334568
$ws.Cells.Item(23, 13).Value = 1.074602483465888
$ws.Cells.Item(23, 14).Value = 1.023703496854225
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.055726307754505
$ws.Cells.Item(24, 4).Value = 1.058158133970626
$ws.Cells.Item(24, 5).Value = 1.068246069016514
$ws.Cells.Item(24, 6).Value = 1.072809863330503
$ws.Cells.Item(24, 9).Value = 1.043169350401046
$ws.Cells.Item(24, 10).Value = 1.061786074811488
$ws.Cells.Item(24, 11).Value = 1.061452649657182
$ws.Cells.Item(24, 12).Value = 1.07150732970224
$ws.Cells.Item(24, 13).Value = 1.076056304521926
$ws.Cells.Item(24, 14).Value = 1.024186394522972
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.057944628691615
$ws.Cells.Item(25, 4).Value = 1.059838234482058
$ws.Cells.Item(25, 5).Value = 1.070337816220624
$ws.Cells.Item(25, 6).Value = 1.074805713342662
$ws.Cells.Item(25, 9).Value = 1.043640240524276
$ws.Cells.Item(25, 10).Value = 1.063407318752291
$ws.Cells.Item(25, 11).Value = 1.062817086506378
$ws.Cells.Item(25, 12).Value = 1.073285617426783
$ws.Cells.Item(25, 13).Value = 1.077740501551991
$ws.Cells.Item(25, 14).Value = 1.024744626748762
